# Resume edit: update the first "MAR 2020 - PRESENTE" / "RESOLVIT INTERNATIONAL"
# job block to "jun 2021 - PRESENTE" / "S1 gateway", and update the second
# "MAR 2020 - PRESENTE" date line to "MAR 2020 – may 2021".
#
# Both blocks are otherwise identical (this resume template repeats the same
# text twice for a layered/overlapping design), so we locate each occurrence
# in document order rather than relying on raw character offsets.

$d = $word.ActiveDocument

function Find-AllRanges($text) {
    $ranges = @()
    $scan = $d.Content
    $scan.Start = 0
    $scan.End = $d.Content.End
    while ($scan.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $ranges += , @($scan.Start, $scan.End)
        $scan.Start = $scan.End
        $scan.End = $d.Content.End
    }
    return $ranges
}

$dateHits = Find-AllRanges("MAR 2020 - PRESENTE")
$titleHits = Find-AllRanges("RESOLVIT INTERNATIONAL")

# First job block: date line immediately followed (a couple paragraphs later)
# by the big "RESOLVIT INTERNATIONAL" company-name paragraph (font size 38
# half-points == 19pt), which is the one that becomes "S1 gateway".
$firstDate = $dateHits[0]
$secondDate = $dateHits[1]

$firstTitle = $null
foreach ($hit in $titleHits) {
    if ($hit[0] -gt $firstDate[1] -and $hit[0] -lt $secondDate[0]) {
        $rng = $d.Range($hit[0], $hit[1])
        if ($rng.Font.Size -eq 19) {
            $firstTitle = $hit
        }
    }
}

# Apply edits from the end of the document backwards so earlier offsets stay valid.

# 3) Second date line: "MAR 2020 - PRESENTE" -> "MAR 2020 – may 2021"
$r3 = $d.Range($secondDate[0], $secondDate[1])
$r3.Find.Execute("MAR 2020 - PRESENTE", $true, $false, $false, $false, $false, $true, 0, $false, "MAR 2020 – may 2021", 2)

# 2) Company name for first block: "RESOLVIT INTERNATIONAL" -> "S1 gateway"
$r2 = $d.Range($firstTitle[0], $firstTitle[1])
$r2.Find.Execute("RESOLVIT INTERNATIONAL", $true, $false, $false, $false, $false, $true, 0, $false, "S1 gateway", 2)

# 1) First date line: "MAR 2020 - PRESENTE" -> "jun 2021 - PRESENTE"
$r1 = $d.Range($firstDate[0], $firstDate[1])
$r1.Find.Execute("MAR 2020 - PRESENTE", $true, $false, $false, $false, $false, $true, 0, $false, "jun 2021 - PRESENTE", 2)

Write-Output "Done."
